# Auto-generated edit script: updates currentAveragePrice / Leve-profit
# columns (H:N) across multiple class sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1887.075
$ws.Range("I137").Value = 1721.6451
$ws.Range("K137").Value = 5164.9353
$ws.Range("M137").Value = -2614.9353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8275.213
$ws.Range("I32").Value = 4650.892
$ws.Range("J32").Value = 31833.3
$ws.Range("K32").Value = 4650.892
$ws.Range("L32").Value = 31833.3
$ws.Range("M32").Value = -4363.892
$ws.Range("N32").Value = -32407.3
$ws.Range("H74").Value = 10139
$ws.Range("I74").Value = 2398.7058
$ws.Range("K74").Value = 2398.7058
$ws.Range("M74").Value = -1524.7058
$ws.Range("H77").Value = 10139
$ws.Range("I77").Value = 2398.7058
$ws.Range("K77").Value = 11993.529
$ws.Range("M77").Value = -7625.529
$ws.Range("H132").Value = 2528.1162
$ws.Range("I132").Value = 2351.081
$ws.Range("J132").Value = 3619.8333
$ws.Range("K132").Value = 7053.243
$ws.Range("L132").Value = 10859.4999
$ws.Range("M132").Value = -4523.243
$ws.Range("N132").Value = -15919.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4527.4546
$ws.Range("I86").Value = 5443.7144
$ws.Range("K86").Value = 5443.7144
$ws.Range("M86").Value = -4320.7144
$ws.Range("H89").Value = 4527.4546
$ws.Range("I89").Value = 5443.7144
$ws.Range("K89").Value = 27218.572
$ws.Range("M89").Value = -21602.572
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4315.6294
$ws.Range("I16").Value = 4063.353
$ws.Range("J16").Value = 4744.5
$ws.Range("K16").Value = 4063.353
$ws.Range("L16").Value = 4744.5
$ws.Range("M16").Value = -3776.353
$ws.Range("N16").Value = -5318.5
$ws.Range("H22").Value = 790.3333
$ws.Range("I22").Value = 226.71428
$ws.Range("J22").Value = 1579.4
$ws.Range("K22").Value = 226.71428
$ws.Range("L22").Value = 1579.4
$ws.Range("M22").Value = 123.28572
$ws.Range("N22").Value = -2279.4
$ws.Range("H31").Value = 1815.5416
$ws.Range("J31").Value = 6969
$ws.Range("L31").Value = 6969
$ws.Range("N31").Value = -7559
$ws.Range("H34").Value = 1815.5416
$ws.Range("J34").Value = 6969
$ws.Range("L34").Value = 6969
$ws.Range("N34").Value = -7373
$ws.Range("H37").Value = 4500
$ws.Range("J37").Value = 4500
$ws.Range("L37").Value = 4500
$ws.Range("N37").Value = -4714
$ws.Range("H86").Value = 17866.727
$ws.Range("I86").Value = 18653.6
$ws.Range("J86").Value = 9998
$ws.Range("K86").Value = 18653.6
$ws.Range("L86").Value = 9998
$ws.Range("M86").Value = -17530.6
$ws.Range("N86").Value = -12244
$ws.Range("H89").Value = 17866.727
$ws.Range("I89").Value = 18653.6
$ws.Range("J89").Value = 9998
$ws.Range("K89").Value = 93268
$ws.Range("L89").Value = 49990
$ws.Range("M89").Value = -87652
$ws.Range("N89").Value = -61222
$ws.Range("H113").Value = 4315.6294
$ws.Range("I113").Value = 4063.353
$ws.Range("J113").Value = 4744.5
$ws.Range("K113").Value = 4063.353
$ws.Range("L113").Value = 4744.5
$ws.Range("M113").Value = -1893.353
$ws.Range("N113").Value = -9084.5
$ws.Range("H127").Value = 100777
$ws.Range("J127").Value = 100777
$ws.Range("L127").Value = 100777
$ws.Range("N127").Value = -110697

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 909188.4
$ws.Range("I7").Value = 113.85714
$ws.Range("J7").Value = 2500068.8
$ws.Range("K7").Value = 341.57142
$ws.Range("L7").Value = 7500206.399999999
$ws.Range("M7").Value = -229.57142
$ws.Range("N7").Value = -7500430.399999999
$ws.Range("H68").Value = 5000911.5
$ws.Range("J68").Value = 5556502
$ws.Range("L68").Value = 16669506
$ws.Range("N68").Value = -16671128
$ws.Range("H71").Value = 5000911.5
$ws.Range("J71").Value = 5556502
$ws.Range("L71").Value = 50008518
$ws.Range("N71").Value = -50016630
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H108").Value = 688
$ws.Range("I108").Value = 688
$ws.Range("K108").Value = 2064
$ws.Range("M108").Value = 816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3737.375
$ws.Range("I122").Value = 999.5
$ws.Range("K122").Value = 2998.5
$ws.Range("M122").Value = -548.5
$ws.Range("H132").Value = 2989.9312
$ws.Range("I132").Value = 2479.5417
$ws.Range("J132").Value = 5439.8
$ws.Range("K132").Value = 7438.625100000001
$ws.Range("L132").Value = 16319.4
$ws.Range("M132").Value = -4908.625100000001
$ws.Range("N132").Value = -21379.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2452.52
$ws.Range("I16").Value = 2931.3157
$ws.Range("K16").Value = 2931.3157
$ws.Range("M16").Value = -2761.3157
$ws.Range("H68").Value = 2758.1428
$ws.Range("I68").Value = 2837.6316
$ws.Range("J68").Value = 2003
$ws.Range("K68").Value = 2837.6316
$ws.Range("L68").Value = 2003
$ws.Range("M68").Value = -2088.6316
$ws.Range("N68").Value = -3501
$ws.Range("H71").Value = 2758.1428
$ws.Range("I71").Value = 2837.6316
$ws.Range("J71").Value = 2003
$ws.Range("K71").Value = 14188.158
$ws.Range("L71").Value = 10015
$ws.Range("M71").Value = -10444.158
$ws.Range("N71").Value = -17503
$ws.Range("H93").Value = 2104.7
$ws.Range("I93").Value = 2108
$ws.Range("J93").Value = 2075
$ws.Range("K93").Value = 2108
$ws.Range("L93").Value = 2075
$ws.Range("M93").Value = -860
$ws.Range("N93").Value = -4571
$ws.Range("H132").Value = 35894.066
$ws.Range("I132").Value = 45601.176
$ws.Range("J132").Value = 3999.2856
$ws.Range("K132").Value = 136803.528
$ws.Range("L132").Value = 11997.8568
$ws.Range("M132").Value = -134273.528
$ws.Range("N132").Value = -17057.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 201799.6
$ws.Range("I107").Value = 2249.5
$ws.Range("J107").Value = 1000000
$ws.Range("K107").Value = 6748.5
$ws.Range("L107").Value = 3000000
$ws.Range("M107").Value = -4828.5
$ws.Range("N107").Value = -3003840
$ws.Range("H132").Value = 1851.7556
$ws.Range("I132").Value = 1746.7778
$ws.Range("K132").Value = 5240.3334
$ws.Range("M132").Value = -2710.3334
